$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.066.38"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").Value = "'3.533.95"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'601.29"
$ws.Range("E5").Value = "  -0.54%  "
$ws.Range("D6").Value = "'145.76"
$ws.Range("E6").Value = "  +0.90%  "
$ws.Range("D7").Value = "'3.533.44"
$ws.Range("E7").Value = "  -0.47%  "
$ws.Range("E8").Value = "  -0.22%  "
$ws.Range("D9").Value = "'0.499"
$ws.Range("E9").Value = "  +1.70%  "
$ws.Range("E10").Value = "  -1.15%  "
$ws.Range("D11").Value = "'7.73"
$ws.Range("E11").Value = "  -1.85%  "
$ws.Range("D12").Value = "'0.406"
$ws.Range("E12").Value = "  -1.45%  "
$ws.Range("D13").Value = "'4.141.62"
$ws.Range("E13").Value = "  -0.28%  "
$ws.Range("D14").Value = "'0.0000201"
$ws.Range("E14").Value = "  -2.84%  "
$ws.Range("D15").Value = "'28.92"
$ws.Range("E15").Value = "  -3.61%  "
$ws.Range("D16").Value = "'3.537.63"
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("E17").Value = "  +1.38%  "
$ws.Range("D18").Value = "'66.077.32"
$ws.Range("E18").Value = "  -0.54%  "
$ws.Range("D19").Value = "'10.99"
$ws.Range("E19").Value = "  -4.47%  "
$ws.Range("D20").Value = "'6.23"
$ws.Range("E20").Value = "  +0.92%  "
$ws.Range("D21").Value = "'14.56"
$ws.Range("E21").Value = "  -1.80%  "
$ws.Range("D22").Value = "'419.56"
$ws.Range("E22").Value = "  -2.45%  "
$ws.Range("D23").Value = "'0.600"
$ws.Range("E23").Value = "  -1.57%  "
$ws.Range("D24").Value = "'77.81"
$ws.Range("E24").Value = "  -2.13%  "
$ws.Range("D25").Value = "'3.677.19"
$ws.Range("E25").Value = "  -0.36%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  -3.26%  "
$ws.Range("D28").Value = "'9.08"
$ws.Range("E28").Value = "  -1.05%  "
$ws.Range("D29").Value = "'2.46"
$ws.Range("E29").Value = "  -1.76%  "
$ws.Range("D30").Value = "'7.74"
$ws.Range("E30").Value = "  -2.72%  "
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("D32").Value = "'3.535.01"
$ws.Range("E32").Value = "  -0.28%  "
$ws.Range("D33").Value = "'0.154"
$ws.Range("E33").Value = "  +0.64%  "
$ws.Range("D34").Value = "'24.29"
$ws.Range("E34").Value = "  -4.50%  "
$ws.Range("D36").Value = "'7.55"
$ws.Range("E36").Value = "  -3.61%  "
$ws.Range("D37").Value = "'1.26"
$ws.Range("E37").Value = "  -12.97%  "
$ws.Range("D38").Value = "'174.26"
$ws.Range("E38").Value = "  -0.93%  "
$ws.Range("D39").Value = "'5.22"
$ws.Range("E39").Value = "  -6.42%  "
$ws.Range("D40").Value = "'1.59"
$ws.Range("E40").Value = "  -7.96%  "
$ws.Range("D41").Value = "'0.0821"
$ws.Range("E41").Value = "  -2.91%  "
$ws.Range("D42").Value = "'5.08"
$ws.Range("E42").Value = "  -1.98%  "
$ws.Range("D43").Value = "'0.858"
$ws.Range("E43").Value = "  -3.37%  "
$ws.Range("D44").Value = "'45.55"
$ws.Range("E44").Value = "  -1.06%  "
$ws.Range("D45").Value = "'1.78"
$ws.Range("E45").Value = "  -6.79%  "
$ws.Range("D46").Value = "'1.00"
$ws.Range("E46").Value = "  +0.14%  "
$ws.Range("E47").Value = "  -4.03%  "
$ws.Range("D48").Value = "'7.09"
$ws.Range("E48").Value = "  -0.54%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'22.53"
$ws.Range("E49").Value = "  -3.39%  "
$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D50").Value = "'1.09"
$ws.Range("E50").Value = "  -8.13%  "
$ws.Range("D51").Value = "'23.04"
$ws.Range("E51").Value = "  -7.87%  "
